$d = $word.ActiveDocument

$d.Content.Find.Execute("39-29=", $true, $false, $false, $false, $false, $true, 1, $false, "40+42=", 2) | Out-Null
$d.Content.Find.Execute("41-36=", $true, $false, $false, $false, $false, $true, 1, $false, "95-7=", 2) | Out-Null
$d.Content.Find.Execute("23+55=", $true, $false, $false, $false, $false, $true, 1, $false, "81-30=", 2) | Out-Null
$d.Content.Find.Execute("26+68=", $true, $false, $false, $false, $false, $true, 1, $false, "70+19=", 2) | Out-Null
$d.Content.Find.Execute("45+30=", $true, $false, $false, $false, $false, $true, 1, $false, "53-6=", 2) | Out-Null
$d.Content.Find.Execute("44+49=", $true, $false, $false, $false, $false, $true, 1, $false, "36+25=", 2) | Out-Null
$d.Content.Find.Execute("20-17=", $true, $false, $false, $false, $false, $true, 1, $false, "85-7=", 2) | Out-Null
$d.Content.Find.Execute("43+41=", $true, $false, $false, $false, $false, $true, 1, $false, "41-10=", 2) | Out-Null
$d.Content.Find.Execute("36+22=", $true, $false, $false, $false, $false, $true, 1, $false, "57-18=", 2) | Out-Null
$d.Content.Find.Execute("45-31=", $true, $false, $false, $false, $false, $true, 1, $false, "89-23=", 2) | Out-Null
$d.Content.Find.Execute("38+2=", $true, $false, $false, $false, $false, $true, 1, $false, "87-6=", 2) | Out-Null
$d.Content.Find.Execute("10+8=", $true, $false, $false, $false, $false, $true, 1, $false, "9+85=", 2) | Out-Null
$d.Content.Find.Execute("88-54=", $true, $false, $false, $false, $false, $true, 1, $false, "9+55=", 2) | Out-Null
$d.Content.Find.Execute("90-20=", $true, $false, $false, $false, $false, $true, 1, $false, "8+9=", 2) | Out-Null
$d.Content.Find.Execute("48+48=", $true, $false, $false, $false, $false, $true, 1, $false, "54-23=", 2) | Out-Null
$d.Content.Find.Execute("89-61=", $true, $false, $false, $false, $false, $true, 1, $false, "90-21=", 2) | Out-Null
$d.Content.Find.Execute("97+2=", $true, $false, $false, $false, $false, $true, 1, $false, "34+59=", 2) | Out-Null
$d.Content.Find.Execute("55+17=", $true, $false, $false, $false, $false, $true, 1, $false, "71-2=", 2) | Out-Null
$d.Content.Find.Execute("58-1=", $true, $false, $false, $false, $false, $true, 1, $false, "32+10=", 2) | Out-Null
$d.Content.Find.Execute("75-44=", $true, $false, $false, $false, $false, $true, 1, $false, "3+79=", 2) | Out-Null
$d.Content.Find.Execute("61-13=", $true, $false, $false, $false, $false, $true, 1, $false, "98-10=", 2) | Out-Null
$d.Content.Find.Execute("53+36=", $true, $false, $false, $false, $false, $true, 1, $false, "8+56=", 2) | Out-Null
$d.Content.Find.Execute("56+30=", $true, $false, $false, $false, $false, $true, 1, $false, "69-56=", 2) | Out-Null
$d.Content.Find.Execute("61-28=", $true, $false, $false, $false, $false, $true, 1, $false, "32+1=", 2) | Out-Null
$d.Content.Find.Execute("47-22=", $true, $false, $false, $false, $false, $true, 1, $false, "5+43=", 2) | Out-Null
$d.Content.Find.Execute("98-74=", $true, $false, $false, $false, $false, $true, 1, $false, "1+69=", 2) | Out-Null
$d.Content.Find.Execute("23+32=", $true, $false, $false, $false, $false, $true, 1, $false, "72-31=", 2) | Out-Null
$d.Content.Find.Execute("34+56=", $true, $false, $false, $false, $false, $true, 1, $false, "38-24=", 2) | Out-Null
$d.Content.Find.Execute("38-8=", $true, $false, $false, $false, $false, $true, 1, $false, "71-36=", 2) | Out-Null
$d.Content.Find.Execute("71-66=", $true, $false, $false, $false, $false, $true, 1, $false, "58-50=", 2) | Out-Null
$d.Content.Find.Execute("87-77=", $true, $false, $false, $false, $false, $true, 1, $false, "95-65=", 2) | Out-Null
$d.Content.Find.Execute("45+27=", $true, $false, $false, $false, $false, $true, 1, $false, "39-22=", 2) | Out-Null
$d.Content.Find.Execute("55+35=", $true, $false, $false, $false, $false, $true, 1, $false, "32+53=", 2) | Out-Null
$d.Content.Find.Execute("56-44=", $true, $false, $false, $false, $false, $true, 1, $false, "56-50=", 2) | Out-Null
$d.Content.Find.Execute("82-74=", $true, $false, $false, $false, $false, $true, 1, $false, "76-24=", 2) | Out-Null
$d.Content.Find.Execute("64-49=", $true, $false, $false, $false, $false, $true, 1, $false, "42+9=", 2) | Out-Null
$d.Content.Find.Execute("31-29=", $true, $false, $false, $false, $false, $true, 1, $false, "4+72=", 2) | Out-Null
$d.Content.Find.Execute("84-83=", $true, $false, $false, $false, $false, $true, 1, $false, "5+52=", 2) | Out-Null
$d.Content.Find.Execute("35+57=", $true, $false, $false, $false, $false, $true, 1, $false, "8-3=", 2) | Out-Null
$d.Content.Find.Execute("76-10=", $true, $false, $false, $false, $false, $true, 1, $false, "53-52=", 2) | Out-Null
$d.Content.Find.Execute("75-60=", $true, $false, $false, $false, $false, $true, 1, $false, "27+2=", 2) | Out-Null
$d.Content.Find.Execute("46+51=", $true, $false, $false, $false, $false, $true, 1, $false, "68-24=", 2) | Out-Null
$d.Content.Find.Execute("10+58=", $true, $false, $false, $false, $false, $true, 1, $false, "55-2=", 2) | Out-Null
$d.Content.Find.Execute("27-13=", $true, $false, $false, $false, $false, $true, 1, $false, "92-11=", 2) | Out-Null
$d.Content.Find.Execute("19+46=", $true, $false, $false, $false, $false, $true, 1, $false, "17+60=", 2) | Out-Null
$d.Content.Find.Execute("89-22=", $true, $false, $false, $false, $false, $true, 1, $false, "79-75=", 2) | Out-Null
$d.Content.Find.Execute("89-13=", $true, $false, $false, $false, $false, $true, 1, $false, "11+21=", 2) | Out-Null
$d.Content.Find.Execute("67-58=", $true, $false, $false, $false, $false, $true, 1, $false, "96-89=", 2) | Out-Null
$d.Content.Find.Execute("61-15=", $true, $false, $false, $false, $false, $true, 1, $false, "96-66=", 2) | Out-Null
$d.Content.Find.Execute("81-77=", $true, $false, $false, $false, $false, $true, 1, $false, "10+80=", 2) | Out-Null
$d.Content.Find.Execute("17+26=", $true, $false, $false, $false, $false, $true, 1, $false, "40+44=", 2) | Out-Null
$d.Content.Find.Execute("47+4=", $true, $false, $false, $false, $false, $true, 1, $false, "39+52=", 2) | Out-Null
$d.Content.Find.Execute("82-75=", $true, $false, $false, $false, $false, $true, 1, $false, "79-41=", 2) | Out-Null
$d.Content.Find.Execute("12-11=", $true, $false, $false, $false, $false, $true, 1, $false, "78-69=", 2) | Out-Null
$d.Content.Find.Execute("42+23=", $true, $false, $false, $false, $false, $true, 1, $false, "89-72=", 2) | Out-Null
$d.Content.Find.Execute("63-1=", $true, $false, $false, $false, $false, $true, 1, $false, "54+4=", 2) | Out-Null
$d.Content.Find.Execute("40+35=", $true, $false, $false, $false, $false, $true, 1, $false, "58-28=", 2) | Out-Null
$d.Content.Find.Execute("15+26=", $true, $false, $false, $false, $false, $true, 1, $false, "18+60=", 2) | Out-Null
$d.Content.Find.Execute("39-39=", $true, $false, $false, $false, $false, $true, 1, $false, "67-32=", 2) | Out-Null
$d.Content.Find.Execute("16-9=", $true, $false, $false, $false, $false, $true, 1, $false, "39+3=", 2) | Out-Null
$d.Content.Find.Execute("90-42=", $true, $false, $false, $false, $false, $true, 1, $false, "42+57=", 2) | Out-Null
$d.Content.Find.Execute("63+16=", $true, $false, $false, $false, $false, $true, 1, $false, "60+7=", 2) | Out-Null
$d.Content.Find.Execute("70-1=", $true, $false, $false, $false, $false, $true, 1, $false, "43-23=", 2) | Out-Null
$d.Content.Find.Execute("23+8=", $true, $false, $false, $false, $false, $true, 1, $false, "38+58=", 2) | Out-Null
$d.Content.Find.Execute("74-50=", $true, $false, $false, $false, $false, $true, 1, $false, "0+49=", 2) | Out-Null
$d.Content.Find.Execute("62-49=", $true, $false, $false, $false, $false, $true, 1, $false, "24+15=", 2) | Out-Null
$d.Content.Find.Execute("14+63=", $true, $false, $false, $false, $false, $true, 1, $false, "34-21=", 2) | Out-Null
$d.Content.Find.Execute("98-58=", $true, $false, $false, $false, $false, $true, 1, $false, "62+24=", 2) | Out-Null
$d.Content.Find.Execute("34+55=", $true, $false, $false, $false, $false, $true, 1, $false, "32+21=", 2) | Out-Null
$d.Content.Find.Execute("33+28=", $true, $false, $false, $false, $false, $true, 1, $false, "79-50=", 2) | Out-Null
$d.Content.Find.Execute("49-2=", $true, $false, $false, $false, $false, $true, 1, $false, "87-24=", 2) | Out-Null
$d.Content.Find.Execute("79-7=", $true, $false, $false, $false, $false, $true, 1, $false, "47-45=", 2) | Out-Null
$d.Content.Find.Execute("40+49=", $true, $false, $false, $false, $false, $true, 1, $false, "27+60=", 2) | Out-Null
$d.Content.Find.Execute("1+97=", $true, $false, $false, $false, $false, $true, 1, $false, "0+66=", 2) | Out-Null
$d.Content.Find.Execute("38+13=", $true, $false, $false, $false, $false, $true, 1, $false, "56-16=", 2) | Out-Null
$d.Content.Find.Execute("7+8=", $true, $false, $false, $false, $false, $true, 1, $false, "12+28=", 2) | Out-Null
$d.Content.Find.Execute("97-7=", $true, $false, $false, $false, $false, $true, 1, $false, "75+14=", 2) | Out-Null
$d.Content.Find.Execute("56-7=", $true, $false, $false, $false, $false, $true, 1, $false, "23+51=", 2) | Out-Null
$d.Content.Find.Execute("66-42=", $true, $false, $false, $false, $false, $true, 1, $false, "75-1=", 2) | Out-Null
$d.Content.Find.Execute("96-90=", $true, $false, $false, $false, $false, $true, 1, $false, "94-48=", 2) | Out-Null
$d.Content.Find.Execute("36-1=", $true, $false, $false, $false, $false, $true, 1, $false, "80-38=", 2) | Out-Null
$d.Content.Find.Execute("1+31=", $true, $false, $false, $false, $false, $true, 1, $false, "82-73=", 2) | Out-Null
$d.Content.Find.Execute("56-53=", $true, $false, $false, $false, $false, $true, 1, $false, "5+54=", 2) | Out-Null
$d.Content.Find.Execute("74-54=", $true, $false, $false, $false, $false, $true, 1, $false, "82-68=", 2) | Out-Null
$d.Content.Find.Execute("42-8=", $true, $false, $false, $false, $false, $true, 1, $false, "76-59=", 2) | Out-Null
$d.Content.Find.Execute("36+0=", $true, $false, $false, $false, $false, $true, 1, $false, "67-18=", 2) | Out-Null
$d.Content.Find.Execute("90-83=", $true, $false, $false, $false, $false, $true, 1, $false, "47+45=", 2) | Out-Null
$d.Content.Find.Execute("1+79=", $true, $false, $false, $false, $false, $true, 1, $false, "77+3=", 2) | Out-Null
$d.Content.Find.Execute("88-35=", $true, $false, $false, $false, $false, $true, 1, $false, "22-13=", 2) | Out-Null
$d.Content.Find.Execute("68+16=", $true, $false, $false, $false, $false, $true, 1, $false, "76+5=", 2) | Out-Null
$d.Content.Find.Execute("58+13=", $true, $false, $false, $false, $false, $true, 1, $false, "69-19=", 2) | Out-Null
$d.Content.Find.Execute("71-31=", $true, $false, $false, $false, $false, $true, 1, $false, "13+71=", 2) | Out-Null
$d.Content.Find.Execute("22+36=", $true, $false, $false, $false, $false, $true, 1, $false, "11+7=", 2) | Out-Null
$d.Content.Find.Execute("82-65=", $true, $false, $false, $false, $false, $true, 1, $false, "35-18=", 2) | Out-Null
$d.Content.Find.Execute("67-30=", $true, $false, $false, $false, $false, $true, 1, $false, "56-49=", 2) | Out-Null
$d.Content.Find.Execute("22-18=", $true, $false, $false, $false, $false, $true, 1, $false, "18-16=", 2) | Out-Null
$d.Content.Find.Execute("89-88=", $true, $false, $false, $false, $false, $true, 1, $false, "49+33=", 2) | Out-Null
$d.Content.Find.Execute("31+29=", $true, $false, $false, $false, $false, $true, 1, $false, "64-19=", 2) | Out-Null
$d.Content.Find.Execute("83+14=", $true, $false, $false, $false, $false, $true, 1, $false, "58+4=", 2) | Out-Null
$d.Content.Find.Execute("52-8=", $true, $false, $false, $false, $false, $true, 1, $false, "54+27=", 2) | Out-Null
